$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.894.56'
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").Value = '1.584.10'
$ws.Range("E3").Value = '  -2.24%  '
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").Value = '''210.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.10%  '
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("D7").Value = '''0.478'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.52%  '
$ws.Range("E8").Value = '  -0.68%  '
$ws.Range("E9").Value = '  -1.62%  '
$ws.Range("D10").Value = '''18.05'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.58%  '
$ws.Range("D11").Value = '''0.0790'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.18%  '
$ws.Range("D12").Value = '1.804.58'
$ws.Range("D13").Value = '1.589.49'
$ws.Range("E13").Value = '  -2.02%  '
$ws.Range("D14").Value = '''4.02'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.89%  '
$ws.Range("E15").Value = '  -2.77%  '
$ws.Range("D16").Value = '25.869.02'
$ws.Range("E16").Value = '  -0.36%  '
$ws.Range("D17").Value = '0.0₃0723'
$ws.Range("E17").Value = '  -1.87%  '
$ws.Range("D18").Value = '''59.76'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.19%  '
$ws.Range("E19").Value = '  -0.30%  '
$ws.Range("D20").Value = '''191.62'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("D22").Value = '''9.34'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.74%  '
$ws.Range("D23").Value = '''5.93'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.54%  '
$ws.Range("E24").Value = '  -0.41%  '
$ws.Range("D25").Value = '''141.23'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.17%  '
$ws.Range("E26").Value = '  -0.26%  '
$ws.Range("E27").Value = '  -0.88%  '
$ws.Range("D28").Value = '''15.06'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.98%  '
$ws.Range("D29").Value = '''6.44'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.24%  '
$ws.Range("E30").Value = '  -5.48%  '
$ws.Range("D31").Value = '''0.0470'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.51%  '
$ws.Range("E32").Value = '  +0.16%  '
$ws.Range("E33").Value = '  -2.27%  '
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("E35").Value = '  -2.41%  '
$ws.Range("D36").Value = '1.096.08'
$ws.Range("E36").Value = '  -2.96%  '
$ws.Range("E37").Value = '  -0.36%  '
$ws.Range("E38").Value = '  -2.28%  '
$ws.Range("E39").Value = '  -1.98%  '
$ws.Range("D40").Value = '''0.500'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.50%  '
$ws.Range("E41").Value = '  -7.80%  '
$ws.Range("D42").Value = '''0.809'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.85%  '
$ws.Range("D43").Value = '''93.51'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.34%  '
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("D45").Value = '1.717.71'
$ws.Range("E45").Value = '  -2.20%  '
$ws.Range("D46").Value = '0.0₆0112'
$ws.Range("E46").Value = '  -0.86%  '
$ws.Range("E47").Value = '  +0.19%  '
$ws.Range("D48").Value = '''53.04'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.96%  '
$ws.Range("E49").Value = '  -1.13%  '
$ws.Range("E50").Value = '  -0.83%  '
$ws.Range("E51").Value = '  -0.27%  '
